# edit.ps1 - applies the hiring-posting text revisions described by the diff:
#  1) "Fischell" gets wrapped in spell-check proofErr markers (run split, no text change)
#  2) The "Proficiency in Information Theory." bullet is expanded to mention
#     Network Science / graph theory.
$d = $word.ActiveDocument

function Get-ParagraphContaining([string]$needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Fischell paragraph: split the run around "Fischell" and mark it with
#    spellStart/spellEnd proofing-error bookmarks, exactly like Word's own
#    spell checker does when it flags a word it doesn't recognise.
# ---------------------------------------------------------------------------
$pFischell = Get-ParagraphContaining "Fischell Department of Bioengineering"
if ($pFischell -eq $null) {
    throw "Could not locate the Fischell paragraph"
}

$fischellXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p w14:paraId="7A040EA4" w14:textId="5236EA82" w:rsidR="007009F3" w:rsidRPr="00830DC7" w:rsidRDefault="009B5251" w:rsidP="007009F3"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="120" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00830DC7"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Imaging- and Neuro-computations for Precision Informatics Research (INSPIRE) Lab is launching at the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Fischell</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Department of Bioengineering at the University of Maryland, College Park!</w:t></w:r><w:r w:rsidR="00481BB0" w:rsidRPr="00830DC7"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r = $pFischell.Range
$r.Text = ""
[void]$r.InsertXML($fischellXml)

# ---------------------------------------------------------------------------
# 2) "Proficiency in Information Theory." bullet: extend it to also mention
#    Network Science / graph theory, split across runs as Word would when a
#    sentence is edited incrementally.
# ---------------------------------------------------------------------------
$pInfoTheory = Get-ParagraphContaining "Proficiency in Information Theory"
if ($pInfoTheory -eq $null) {
    throw "Could not locate the Proficiency in Information Theory paragraph"
}

$infoTheoryXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p w14:paraId="74056F5C" w14:textId="1D01A009" w:rsidR="001D412A" w:rsidRDefault="001D412A" w:rsidP="001D412A"><w:pPr><w:widowControl/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:autoSpaceDE/><w:autoSpaceDN/><w:ind w:left="360"/><w:jc w:val="both"/><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Proficiency in Information Theory</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> and/or</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Network Science</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> (e.g., graph theory).</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r2 = $pInfoTheory.Range
$r2.Text = ""
[void]$r2.InsertXML($infoTheoryXml)

Write-Output "Applied Fischell proofing split and Information Theory bullet expansion."
